# Update countries & provincias Spain
# Refresh the COVID country table: update the "last updated" timestamp
# and re-sync the rows whose ranking changed after refreshing the
# underlying case counts (ties/overtakes reorder some country rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 19:35"

# Rows that changed: row number, country name, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(4,   "Estados Unidos",  6946888, 20947, 4201211, 2542194, 0, 312, 203483),
    @(5,   "India",           5367361, 61886, 4265768, 1015388, 0, 580,  86205),
    @(25,  "Alemania",         271840,   596,  243000,   19374, 0,   2,   9466),
    @(38,  "Marruecos",         99816,  2552,   79008,   19013, 0,  40,   1795),
    @(39,  "Belgica",           99649,  1673,   18908,   70804, 0,   1,   9937),
    @(40,  "Kuwait",            99049,   521,   89498,    8970, 0,   1,    581),
    @(51,  "Etiopia",           68131,   616,   27939,   39103, 0,  17,   1089),
    @(52,  "Portugal",          68025,   849,   45404,   20722, 0,   5,   1899),
    @(59,  "Uzbekistan",        50872,   619,   47121,    3326, 0,   6,    425),
    @(60,  "Argelia",           49623,   210,   34923,   13035, 0,   6,   1665),
    @(62,  "Chequia",           47285,  1023,   24018,   22768, 0,   4,    499),
    @(63,  "Armenia",           47154,   244,   42551,    3675, 0,   2,    928),
    @(64,  "Moldavia",          46336,   688,   34236,   10899, 0,  15,   1201),
    @(73,  "Irlanda",           32538,   267,   23364,    7382, 0,   0,   1792),
    @(75,  "Libano",            28297,   779,   11440,   16571, 0,   5,    286),
    @(111, "Mozambique",         6537,   273,    3620,    2876, 0,   1,     41),
    @(134, "Gambia",              3504,   19,    1992,    1404, 0,   0,    108),
    @(135, "Tailandia",           3500,    3,    3338,     103, 0,   1,     59),
    @(140, "Sri Lanka",           3283,    2,    3070,     200, 0,   0,     13),
    @(153, "Yemen",               2026,    2,    1221,     220, 0,   0,    585),
    @(204, "Timor Oriental",        27,    0,      26,       1, 0,   0,      0),
    @(205, "Santa Lucia",           27,    0,      26,       1, 0,   0,      0),
    @(214, "Islas Malvinas",        13,    0,      13,       0, 0,   0,      0),
    @(215, "Montserrat",            13,    0,      12,       0, 0,   0,      1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]
    $ws.Range("H$r").Value = $row[8]
}
